# Updates the running-results sheet: refreshed activity data (dates, times,
# locations, distances), table autofilter, uniform column widths, and
# removal of the old conditional-formatting rules.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-row data: row -> (Start/End-Date+Time serial #1, serial #2, Location, Distance)
$data = @{}
$data[2] = @(43101.4108352662, 43101.43955054398, "Antwerp", 8.516109537338451)
$data[3] = @(43107.17827626158, 43107.20160959491, "Heusden-Zolder", 7.347180371568261)
$data[4] = @(43108.245410034724, 43108.26520170139, "Ghent", 5.021017254063722)
$data[5] = @(43117.001220729166, 43117.021243877316, "Schulen", 5.483717788175882)
$data[6] = @(43117.09622403935, 43117.14697635417, "Lummen", 12.198801960702552)
$data[7] = @(43125.788834189814, 43125.806334189816, "Schulen", 5.432378764660736)
$data[8] = @(43134.27033746528, 43134.30185366898, "Ghent", 10.048111582066234)
$data[9] = @(43148.45178474537, 43148.49251391204, "Heusden-Zolder", 13.0071459362053)
$data[10] = @(43157.8351708912, 43157.86160607639, "Schulen", 7.150295247138457)
$data[11] = @(43161.91693673611, 43161.951716828706, "Heusden-Zolder", 9.811889009469633)
$data[12] = @(43162.55858101852, 43162.59043287037, "Antwerp", 10.078978434153257)
$data[13] = @(43163.06331835648, 43163.09145493055, "Lummen", 9.111168021047215)
$data[14] = @(43171.775890046294, 43171.799686342594, "Antwerp", 7.639479810321195)
$data[15] = @(43174.887463541665, 43174.90484780093, "Antwerp", 5.102015479309499)
$data[16] = @(43178.86926556713, 43178.89737899305, "Antwerp", 6.701229977060974)
$data[17] = @(43181.19306818287, 43181.235487164355, "Brussels", 12.59215310002132)
$data[18] = @(43181.70130883102, 43181.73355420139, "Heusden-Zolder", 8.432484678444226)
$data[19] = @(43186.534220092595, 43186.56221777778, "Brussels", 7.902374773276685)
$data[20] = @(43186.99159445602, 43187.04036760417, "Ghent", 11.554691117278642)
$data[21] = @(43196.38063280092, 43196.402322615744, "Schulen", 5.2430426184446315)
$data[22] = @(43200.46987519676, 43200.48728260417, "Brussels", 5.005267193320042)
$data[23] = @(43206.94643818287, 43206.99643818287, "Heusden-Zolder", 12.52107801286148)
$data[24] = @(43209.32189533565, 43209.353585150464, "Brussels", 9.831140704045616)
$data[25] = @(43209.38638125, 43209.40896226852, "Ghent", 6.996415265612087)
$data[26] = @(43215.25432825232, 43215.29827501158, "Antwerp", 11.444874747463466)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]   # Start Date
    $ws.Cells.Item($r, 2).Value = $vals[0]   # Start Time
    $ws.Cells.Item($r, 3).Value = $vals[1]   # End Date
    $ws.Cells.Item($r, 4).Value = $vals[1]   # End Time
    $ws.Cells.Item($r, 5).Value = $vals[2]   # Location
    $ws.Cells.Item($r, 6).Value = $vals[3]   # Distance
}

# Turn on the table's autofilter dropdowns
$lo = $ws.ListObjects.Item(1)
$lo.ShowAutoFilter = $true

# Make every column a uniform 12.0 characters wide (drop the old bestFit widths).
# ColumnWidth is offset by 5/6 from the stored sheet width, so 12 - 5/6 round-trips to 12.0.
for ($c = 1; $c -le 8; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 11.166666666666666
}

# Drop the three old conditional-formatting rules (color scale / icon set / data bar).
$ws.Range("H2:H27").FormatConditions.Delete()
$ws.Range("F2:F27").FormatConditions.Delete()
$ws.Range("G2:G27").FormatConditions.Delete()
